# Apply "Trade #52 closed" update to the live trading results workbook.
# - Updates aggregate statistics on the "Summary" sheet.
# - Updates the MarketMaking strategy row on the "Strategy Status" sheet.
# - Appends the new trade row to both the "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.52   # Current Capital
$summary.Range("B4").Value = 0.52      # Total P&L $
$summary.Range("B5").Value = 0.2       # Total P&L %
$summary.Range("B6").Value = 52        # Total Trades
$summary.Range("B8").Value = 28        # Losing Trades
$summary.Range("B9").Value = 28.85     # Win Rate %

# ---------------------------------------------------------------------------
# 2. Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.52     # Capital
$status.Range("D4").Value = 52         # Trades
$status.Range("E4").Value = 0.52       # P&L $
$status.Range("F4").Value = 0.52       # P&L %
$status.Range("G4").Value = 28.85      # Win Rate %

# ---------------------------------------------------------------------------
# 3. Append new trade row (row 53) to "All Trades" and "MarketMaking" sheets
# ---------------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 53

    $ws.Range("A$row").Value = 52

    # Column B holds a date stored as plain text in this workbook (not a
    # native Excel date value), so force a text format before assignment to
    # stop Excel from auto-converting it to a date serial, then clear the
    # formatting again so the cell keeps the sheet's default (unstyled) look.
    $ws.Range("B$row").NumberFormat = "@"
    $ws.Range("B$row").Value = "2026-02-17"
    $ws.Range("B$row").ClearFormats()
    $ws.Range("C$row").Value = "15:42:46"

    $ws.Range("D$row").Value = "MarketMaking"
    $ws.Range("E$row").Value = "DOWN"
    $ws.Range("F$row").Value = 0.68
    $ws.Range("G$row").Value = 0.66
    $ws.Range("H$row").Value = "CLOSED"
    $ws.Range("I$row").Value = -2.9412
    $ws.Range("J$row").Value = -0.02
    $ws.Range("K$row").Value = 100.52
    $ws.Range("L$row").Value = 0
    $ws.Range("M$row").Value = 0
    $ws.Range("N$row").Value = 0.6
    $ws.Range("O$row").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P$row").Value = "early_exit"
    $ws.Range("Q$row").Value = 0.15
}
